# Add results and unfolding with 100 keV threshold
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Beta" row (row 2) values that changed due to the
# re-run of the unfolding with the new 100 keV threshold.
$ws.Range("C2").Value = 19.12075701903682
$ws.Range("E2").Value = 0.01982943797740053
$ws.Range("F2").Value = 10.58302336205827
$ws.Range("G2").Value = 10.1852961245796
$ws.Range("H2").Value = 10.98359017836801
$ws.Range("I2").Value = 0.002323829143158643
$ws.Range("J2").Value = 0.00119761819857747
$ws.Range("K2").Value = 0.003777740334513748
$ws.Range("L2").Value = 0.01054190307816003
$ws.Range("M2").Value = 0.009940987560958962
$ws.Range("N2").Value = 0.01117208578707817

# Update existing "Gamma" row (row 3) values
$ws.Range("C3").Value = 0.04981522627320694
$ws.Range("D3").Value = 0.04815098319456564
$ws.Range("E3").Value = 0.0499839736740351
$ws.Range("F3").Value = 0.04772774741515929
$ws.Range("G3").Value = 0.04740074944560187
$ws.Range("H3").Value = 0.0480437614667096
$ws.Range("I3").Value = 0.04613160149590845
$ws.Range("J3").Value = 0.04581664879518391
$ws.Range("K3").Value = 0.04643484175287912
$ws.Range("L3").Value = 0.04775963952245733
$ws.Range("M3").Value = 0.04743270697781536
$ws.Range("N3").Value = 0.04807564555405727

# Add new row 4 for "Beta + Gamma" combined results
$ws.Range("A4").Value = 2
# Match the formatting used by the existing "Particle Type Id" cells
# (bold, centered, bordered) by copying A2's format onto A4.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 19.17057224531003
$ws.Range("D4").Value = 0.05509537454402212
$ws.Range("E4").Value = 0.06981341165143562
$ws.Range("F4").Value = 10.63075110947343
$ws.Range("G4").Value = 10.2326968740252
$ws.Range("H4").Value = 11.03163393983472
$ws.Range("I4").Value = 0.0484554306390671
$ws.Range("J4").Value = 0.04701426699376137
$ws.Range("K4").Value = 0.05021258208739286
$ws.Range("L4").Value = 0.05830154260061737
$ws.Range("M4").Value = 0.05737369453877433
$ws.Range("N4").Value = 0.05924773134113544
